{"js": "// Update the two-digit-division answer table: each cell's \"a\u00f7b=c, d\"\n// text is replaced by a new equation, in document order. Several old\n// values repeat (e.g. \"14\u00f73=4, 2\" appears twice), so replacements are\n// applied strictly in document order: for each mapping entry we search\n// for the old text and replace the first remaining match, which makes\n// subsequent duplicate matches line up with later mapping entries.\nconst replacements = [\n  [\"29\u00f78=3, 5\", \"90\u00f75=18, 0\"],\n  [\"18\u00f73=6, 0\", \"81\u00f79=9, 0\"],\n  [\"33\u00f73=11, 0\", \"88\u00f79=9, 7\"],\n  [\"40\u00f76=6, 4\", \"56\u00f74=14, 0\"],\n  [\"23\u00f78=2, 7\", \"22\u00f73=7, 1\"],\n  [\"45\u00f77=6, 3\", \"47\u00f77=6, 5\"],\n  [\"14\u00f73=4, 2\", \"42\u00f74=10, 2\"],\n  [\"14\u00f73=4, 2\", \"25\u00f77=3, 4\"],\n  [\"55\u00f75=11, 0\", \"79\u00f74=19, 3\"],\n  [\"31\u00f76=5, 1\", \"81\u00f75=16, 1\"],\n  [\"83\u00f79=9, 2\", \"66\u00f75=13, 1\"],\n  [\"55\u00f74=13, 3\", \"92\u00f77=13, 1\"],\n  [\"88\u00f72=44, 0\", \"91\u00f75=18, 1\"],\n  [\"48\u00f72=24, 0\", \"62\u00f79=6, 8\"],\n  [\"58\u00f72=29, 0\", \"39\u00f75=7, 4\"],\n  [\"78\u00f75=15, 3\", \"98\u00f73=32, 2\"],\n  [\"95\u00f73=31, 2\", \"48\u00f79=5, 3\"],\n  [\"57\u00f73=19, 0\", \"60\u00f78=7, 4\"],\n  [\"66\u00f73=22, 0\", \"31\u00f74=7, 3\"],\n  [\"42\u00f77=6, 0\", \"63\u00f79=7, 0\"],\n  [\"96\u00f72=48, 0\", \"98\u00f75=19, 3\"],\n  [\"58\u00f72=29, 0\", \"78\u00f76=13, 0\"],\n  [\"87\u00f74=21, 3\", \"59\u00f79=6, 5\"],\n  [\"78\u00f75=15, 3\", \"55\u00f77=7, 6\"],\n  [\"81\u00f72=40, 1\", \"52\u00f72=26, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (!results.items.length) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  // Always take the first remaining occurrence - earlier mapping entries\n  // consume earlier-in-document duplicates first, so this keeps every\n  // replacement aligned with the correct table cell.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-division answer table: each cell's \"a\u00f7b=c, d\"\n# text is replaced by a new equation, in document order. Several old\n# values repeat (e.g. \"14\u00f73=4, 2\" appears twice), so each mapping entry\n# is applied with wdReplaceOne (1) against a freshly-fetched Content\n# range, which matches/replaces only the first remaining occurrence -\n# this keeps every replacement lined up with the correct table cell\n# even when the search text is not unique.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"29\u00f78=3, 5\", \"90\u00f75=18, 0\"),\n    @(\"18\u00f73=6, 0\", \"81\u00f79=9, 0\"),\n    @(\"33\u00f73=11, 0\", \"88\u00f79=9, 7\"),\n    @(\"40\u00f76=6, 4\", \"56\u00f74=14, 0\"),\n    @(\"23\u00f78=2, 7\", \"22\u00f73=7, 1\"),\n    @(\"45\u00f77=6, 3\", \"47\u00f77=6, 5\"),\n    @(\"14\u00f73=4, 2\", \"42\u00f74=10, 2\"),\n    @(\"14\u00f73=4, 2\", \"25\u00f77=3, 4\"),\n    @(\"55\u00f75=11, 0\", \"79\u00f74=19, 3\"),\n    @(\"31\u00f76=5, 1\", \"81\u00f75=16, 1\"),\n    @(\"83\u00f79=9, 2\", \"66\u00f75=13, 1\"),\n    @(\"55\u00f74=13, 3\", \"92\u00f77=13, 1\"),\n    @(\"88\u00f72=44, 0\", \"91\u00f75=18, 1\"),\n    @(\"48\u00f72=24, 0\", \"62\u00f79=6, 8\"),\n    @(\"58\u00f72=29, 0\", \"39\u00f75=7, 4\"),\n    @(\"78\u00f75=15, 3\", \"98\u00f73=32, 2\"),\n    @(\"95\u00f73=31, 2\", \"48\u00f79=5, 3\"),\n    @(\"57\u00f73=19, 0\", \"60\u00f78=7, 4\"),\n    @(\"66\u00f73=22, 0\", \"31\u00f74=7, 3\"),\n    @(\"42\u00f77=6, 0\", \"63\u00f79=7, 0\"),\n    @(\"96\u00f72=48, 0\", \"98\u00f75=19, 3\"),\n    @(\"58\u00f72=29, 0\", \"78\u00f76=13, 0\"),\n    @(\"87\u00f74=21, 3\", \"59\u00f79=6, 5\"),\n    @(\"78\u00f75=15, 3\", \"55\u00f77=7, 6\"),\n    @(\"81\u00f72=40, 1\", \"52\u00f72=26, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: '$oldText'\"\n    }\n}\n"}
